$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated rf-adjusted statistics (row => B,D,E,F values)
$data = @{
    2 = @(11.55418027109837, 0.006844427974319675, 1.159950107606476, 1.345484252136276)
    3 = @(11.5575747443692, 0.006833404722591967, 1.158081957035573, 1.341153819211343)
    4 = @(11.43594773668599, 0.006859220858480134, 1.162457111499002, 1.351306536074604)
    5 = @(11.35588598269685, 0.007117919130883461, 1.206299648821002, 1.455158842745673)
    6 = @(11.43882831403803, 0.006850669970385716, 1.161007961387048, 1.347939486404109)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("D$row").Value = $vals[1]
    $ws.Range("E$row").Value = $vals[2]
    $ws.Range("F$row").Value = $vals[3]
}
